$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update price (D) and volume-change (E) cells for rows whose coin identity is unchanged ---
$ws.Range("D2").Value = "68.101.50"
$ws.Range("E2").Value = "  +1.40%  "

$ws.Range("D3").Value = "3.533.71"
$ws.Range("E3").Value = "  +0.40%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'601.76"
$ws.Range("E5").Value = "  +1.43%  "

$ws.Range("D6").Value = "'184.55"
$ws.Range("E6").Value = "  +6.20%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'0.597"
$ws.Range("E8").Value = "  +0.50%  "

$ws.Range("E9").Value = "  +4.31%  "

$ws.Range("E10").Value = "  -1.18%  "

$ws.Range("E11").Value = "  +1.93%  "

$ws.Range("D12").Value = "4.142.79"
$ws.Range("E12").Value = "  +0.29%  "

$ws.Range("D13").Value = "'32.59"
$ws.Range("E13").Value = "  +11.93%  "

$ws.Range("E14").Value = "  -0.25%  "

$ws.Range("D15").Value = "67.998.64"
$ws.Range("E15").Value = "  +1.26%  "

$ws.Range("E16").Value = "  +0.68%  "

$ws.Range("D17").Value = "3.537.67"
$ws.Range("E17").Value = "  -0.74%  "

$ws.Range("D18").Value = "'6.42"
$ws.Range("E18").Value = "  +1.31%  "

$ws.Range("D19").Value = "'14.80"
$ws.Range("E19").Value = "  +3.78%  "

$ws.Range("D20").Value = "'401.07"
$ws.Range("E20").Value = "  +1.33%  "

$ws.Range("D21").Value = "'8.13"
$ws.Range("E21").Value = "  +1.39%  "

$ws.Range("D22").Value = "'73.93"
$ws.Range("E22").Value = "  +1.22%  "

$ws.Range("D23").Value = "'0.548"
$ws.Range("E23").Value = "  +1.17%  "

$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  -0.04%  "

$ws.Range("D25").Value = "'5.68"
$ws.Range("E25").Value = "  -0.40%  "

$ws.Range("E26").Value = "  +2.34%  "

$ws.Range("D27").Value = "'10.72"
$ws.Range("E27").Value = "  +3.84%  "

$ws.Range("E28").Value = "  -0.92%  "

$ws.Range("E29").Value = "  -0.25%  "

$ws.Range("E32").Value = "  +1.15%  "

$ws.Range("D33").Value = "'24.17"
$ws.Range("E33").Value = "  +0.98%  "

$ws.Range("D34").Value = "'7.48"
$ws.Range("E34").Value = "  +1.58%  "

$ws.Range("E35").Value = "  +0.09%  "

$ws.Range("E36").Value = "  +1.86%  "

$ws.Range("D37").Value = "'164.35"
$ws.Range("E37").Value = "  +0.75%  "

$ws.Range("D38").Value = "'0.885"
$ws.Range("E38").Value = "  -1.50%  "

$ws.Range("E39").Value = "  +2.15%  "

$ws.Range("D40").Value = "'7.17"
$ws.Range("E40").Value = "  +2.39%  "

$ws.Range("D41").Value = "'2.83"
$ws.Range("E41").Value = "  +6.95%  "

$ws.Range("D44").Value = "2.887.09"

$ws.Range("D47").Value = "'42.57"
$ws.Range("E47").Value = "  -0.91%  "

$ws.Range("D50").Value = "'1.09"
$ws.Range("E50").Value = "  -0.41%  "

$ws.Range("D51").Value = "'34.12"
$ws.Range("E51").Value = "  +1.75%  "

# --- Rows whose coin identity (B/C) changed position, plus new D/E values ---
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'6.34"
$ws.Range("E30").Value = "  +0.78%  "

$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "'1.49"
$ws.Range("E31").Value = "  +1.98%  "

$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "'27.32"
$ws.Range("E42").Value = "  +3.14%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'4.77"
$ws.Range("E43").Value = "  +1.31%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'27.51"
$ws.Range("E45").Value = "  -0.39%  "

$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "'0.0747"
$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0308"
$ws.Range("E48").Value = "  +0.77%  "

$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").Value = "'349.87"
$ws.Range("E49").Value = "  +3.77%  "
